# Updates the crypto price/volume snapshot (GitHub Actions refresh).
# Numeric-looking "Price" values are written with a leading "'" so Excel's
# COM layer stores them as text (matching the original inlineStr cells)
# instead of coercing them to numbers; Style is reset to "Normal" right
# after so no stray quote-prefix style id is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.934.87"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "2.316.03"
$ws.Range("E3").Value = "  -4.12%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'549.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'131.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("D8").Value = "'0.573"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.68%  "
$ws.Range("D9").Value = "2.315.15"
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("D10").Value = "'0.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").Value = "'5.56"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  -5.04%  "
$ws.Range("D14").Value = "'24.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "2.727.56"
$ws.Range("E15").Value = "  -4.25%  "
$ws.Range("D16").Value = "58.863.34"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("E17").Value = "  -2.92%  "
$ws.Range("D18").Value = "2.311.31"
$ws.Range("E18").Value = "  -4.70%  "
$ws.Range("D19").Value = "'10.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.56%  "
$ws.Range("E20").Value = "  -3.87%  "
$ws.Range("D21").Value = "'315.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("E22").Value = "  -4.20%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'63.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("D25").Value = "'0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.55%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "'8.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.16%  "
$ws.Range("E28").Value = "  -6.76%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'169.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E32").Value = "  +3.85%  "
$ws.Range("D33").Value = "'5.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.75%  "
$ws.Range("D34").Value = "'0.385"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D36").Value = "'17.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("E38").Value = "  -4.69%  "
$ws.Range("D39").Value = "'4.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.55%  "
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("D42").Value = "'304.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.28%  "
$ws.Range("D43").Value = "'141.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("E44").Value = "  -5.33%  "
$ws.Range("D45").Value = "'0.0952"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "'0.0503"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.561"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'18.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.65%  "
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("D50").Value = "'16.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.21%  "
$ws.Range("D51").Value = "'11.03"
$ws.Range("D51").Style = "Normal"
